$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("G6").Value = 1.75
$ws.Range("H6").Value = 3.25
$ws.Range("I6").Value = 4.85
$ws.Range("L6").Value = 1.39
$ws.Range("M6").Value = 2.57
$ws.Range("N6").Value = 2.12
$ws.Range("O6").Value = 1.57
$ws.Range("P6").Value = 1.44
$ws.Range("Q6").Value = 2.4
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 1.65
$ws.Range("T6").Value = 5.5
$ws.Range("U6").Value = 7.1
$ws.Range("W6").Value = 13.5
$ws.Range("Z6").Value = 7.6
$ws.Range("AB6").Value = 18
$ws.Range("AC6").Value = 110
$ws.Range("AE6").Value = 11.25
$ws.Range("AF6").Value = 27
$ws.Range("AG6").Value = 16
$ws.Range("AH6").Value = 90
$ws.Range("AI6").Value = 55
$ws.Range("AJ6").Value = 60

# Row 7
$ws.Range("G7").Value = 1.47
$ws.Range("H7").Value = 3.95
$ws.Range("I7").Value = 6.2
$ws.Range("L7").Value = 1.33
$ws.Range("M7").Value = 2.8
$ws.Range("N7").Value = 1.98
$ws.Range("O7").Value = 1.65
$ws.Range("P7").Value = 1.38
$ws.Range("Q7").Value = 2.6
$ws.Range("R7").Value = 2.15
$ws.Range("S7").Value = 1.55
$ws.Range("T7").Value = 5.5
$ws.Range("U7").Value = 6
$ws.Range("V7").Value = 8.75
$ws.Range("W7").Value = 9.5
$ws.Range("X7").Value = 13.5
$ws.Range("Y7").Value = 37
$ws.Range("Z7").Value = 8.75
$ws.Range("AA7").Value = 7.9
$ws.Range("AB7").Value = 23
$ws.Range("AC7").Value = 150
$ws.Range("AE7").Value = 13.5
$ws.Range("AF7").Value = 37
$ws.Range("AG7").Value = 21
$ws.Range("AH7").Value = 150
$ws.Range("AI7").Value = 80
$ws.Range("AJ7").Value = 90

# Row 10
$ws.Range("L10").Value = 1.25
$ws.Range("M10").Value = 3.75
$ws.Range("N10").Value = 1.85
$ws.Range("O10").Value = 2

# Row 11
$ws.Range("AG11").Value = 12

# Row 13
$ws.Range("J13").Value = 1.02
$ws.Range("L13").Value = 1.13

# Row 15
$ws.Range("G15").Value = 2.27
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = 3.4
$ws.Range("J15").Value = 1.11
$ws.Range("K15").Value = 5.9
$ws.Range("L15").Value = 1.47
$ws.Range("M15").Value = 2.55
$ws.Range("N15").Value = 2.4
$ws.Range("O15").Value = 1.53
$ws.Range("P15").Value = 1.53
$ws.Range("Q15").Value = 2.42
$ws.Range("R15").Value = 2
$ws.Range("S15").Value = 1.72
$ws.Range("T15").Value = 6
$ws.Range("U15").Value = 10.25
$ws.Range("V15").Value = 10
$ws.Range("W15").Value = 25
$ws.Range("X15").Value = 24
$ws.Range("Y15").Value = 45
$ws.Range("Z15").Value = 5.9
$ws.Range("AA15").Value = 6.1
$ws.Range("AB15").Value = 18
$ws.Range("AC15").Value = 120
$ws.Range("AE15").Value = 8.25
$ws.Range("AG15").Value = 13
$ws.Range("AH15").Value = 55
$ws.Range("AI15").Value = 40
$ws.Range("AJ15").Value = 55

# Row 19
$ws.Range("G19").Value = 2.22
$ws.Range("I19").Value = 3.15
$ws.Range("L19").Value = 1.35
$ws.Range("M19").Value = 2.72
$ws.Range("N19").Value = 2.02
$ws.Range("O19").Value = 1.62
$ws.Range("P19").Value = 1.45
$ws.Range("Q19").Value = 2.37
$ws.Range("R19").Value = 1.78
$ws.Range("S19").Value = 1.82
$ws.Range("T19").Value = 6.9
$ws.Range("U19").Value = 10.25
$ws.Range("V19").Value = 9
$ws.Range("W19").Value = 22
$ws.Range("X19").Value = 19.5
$ws.Range("Z19").Value = 8.25
$ws.Range("AC19").Value = 75
$ws.Range("AE19").Value = 9
$ws.Range("AF19").Value = 16
$ws.Range("AG19").Value = 11
$ws.Range("AH19").Value = 40
$ws.Range("AI19").Value = 28
$ws.Range("AJ19").Value = 37

# Row 20
$ws.Range("G20").Value = 1.06
$ws.Range("H20").Value = 13
$ws.Range("I20").Value = 29
$ws.Range("R20").Value = 1.91
$ws.Range("S20").Value = 1.8
$ws.Range("U20").Value = 11
$ws.Range("W20").Value = 8.5
$ws.Range("X20").Value = 12
$ws.Range("Y20").Value = 29
$ws.Range("AA20").Value = 26
$ws.Range("AB20").Value = 34
$ws.Range("AC20").Value = 67
$ws.Range("AF20").Value = 151
$ws.Range("AI20").Value = 151
$ws.Range("AJ20").Value = 101

# Row 23
$ws.Range("G23").Value = 3.6
$ws.Range("H23").Value = 3.6
$ws.Range("M23").Value = 3.9
$ws.Range("P23").Value = 1.34
$ws.Range("Q23").Value = 3
$ws.Range("T23").Value = 13.5
$ws.Range("W23").Value = 50
$ws.Range("X23").Value = 28
$ws.Range("AA23").Value = 7.1
$ws.Range("AH23").Value = 17
